# "Created Account test added"
# The workbook previously had a "searchSomething" test sheet that is no
# longer needed, and the "CreateAccountTest" sheet is renamed to the
# camelCase "createAccountTest" naming convention used by the other sheets.

$wb = $excel.ActiveWorkbook

# Remove the obsolete "searchSomething" worksheet entirely (its only
# shared strings - "key"/"alak" - go away with it).
$wb.Worksheets("searchSomething").Delete() | Out-Null

# Rename "CreateAccountTest" -> "createAccountTest" and make it the
# active / selected sheet with the cursor back at the top-left cell.
$ws = $wb.Worksheets("CreateAccountTest")
$ws.Name = "createAccountTest"
$ws.Activate() | Out-Null
$ws.Range("A1").Select() | Out-Null
